$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row additions/shifts (row 1 unchanged) ---

# Row 2-5: add word_type ("generic") in column J, existing K values stay same text
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# Row 6-13: column C becomes "generic" (word_type), column D becomes "can"/"do"/"look"/"where" (carrier)
$ws.Range("C6").Value = "generic"
$ws.Range("D6").Value = "can"

$ws.Range("C7").Value = "generic"
$ws.Range("D7").Value = "can"

$ws.Range("C8").Value = "generic"
$ws.Range("D8").Value = "do"

$ws.Range("C9").Value = "generic"
$ws.Range("D9").Value = "do"

$ws.Range("C10").Value = "generic"
$ws.Range("D10").Value = "look"

$ws.Range("C11").Value = "generic"
$ws.Range("D11").Value = "look"

$ws.Range("C12").Value = "generic"
$ws.Range("D12").Value = "where"

$ws.Range("C13").Value = "generic"
$ws.Range("D13").Value = "where"

# --- Add new block at the bottom (rows 27-36) ---

$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
